# Update "want to go" counts (column F) for several events.
# Sheet "展览" (Exhibition): F3, F11, F12 change, F2 unchanged.
# Sheet "全部类型" (All Types): F2, F3, F11, F12 all change.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 596
$wsExhibition.Range("F11").Value = 4792
$wsExhibition.Range("F12").Value = 4533

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 730
$wsAll.Range("F3").Value = 596
$wsAll.Range("F11").Value = 4792
$wsAll.Range("F12").Value = 4533
